$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move header "upsa_email" from G1 to E1
$ws.Range("E1").Value2 = $ws.Range("G1").Value2
$ws.Range("G1").ClearContents()

# Move value "mary@epam.com" from G3 to E3
$ws.Range("E3").Value2 = $ws.Range("G3").Value2
$ws.Range("G3").ClearContents()

# Adjust column widths: column E should take the bestFit width that column G had
# (~16.71 OOXML width units -> ~15.83 in the COM ColumnWidth scale), and column G
# should return to the workbook's standard (default) column width.
$ws.Columns("E").ColumnWidth = 15.83
$ws.Columns("G").ColumnWidth = $ws.StandardWidth

# Update the selected cell/range
$ws.Range("G6").Select() | Out-Null
